# Update marksheet correct/total marks figures on the "quiz" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: number of right-answer points per question changed 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total marks scored changed 66 -> 110, and the
# "scored/out-of" label updated from "66/84" to "110/140"
$ws.Range("B12").Value = 110
$ws.Range("E12").Value = "110/140"
